# Append 4 new package rows (94-97) to the "Packages" worksheet,
# mirroring the rows already present in the sheet (columns A-G):
#   A: PackageType   B: State   C: Authority   D: ActionType
#   E: PackageID     F: Status  G: ParentID

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 94 - SPA / MD / Medicaid SPA / (blank) / MD-25-9565 / (blank) / (blank)
$ws.Cells.Item(94, 1).Value = "SPA"
$ws.Cells.Item(94, 2).Value = "MD"
$ws.Cells.Item(94, 3).Value = "Medicaid SPA"
$ws.Cells.Item(94, 5).Value = "MD-25-9565"

# Row 95 - SPA / MD / CHIP SPA / (blank) / MD-25-9566 / (blank) / (blank)
$ws.Cells.Item(95, 1).Value = "SPA"
$ws.Cells.Item(95, 2).Value = "MD"
$ws.Cells.Item(95, 3).Value = "CHIP SPA"
$ws.Cells.Item(95, 5).Value = "MD-25-9566"

# Row 96 - Waiver / MD / 1915(c) / Amendment / MD-2260.R00.73 / (blank) / MD-2260.R00.00
$ws.Cells.Item(96, 1).Value = "Waiver"
$ws.Cells.Item(96, 2).Value = "MD"
$ws.Cells.Item(96, 3).Value = "1915(c)"
$ws.Cells.Item(96, 4).Value = "Amendment"
$ws.Cells.Item(96, 5).Value = "MD-2260.R00.73"
$ws.Cells.Item(96, 7).Value = "MD-2260.R00.00"

# Row 97 - Waiver / MD / 1915(c) / Amendment / MD-2260.R00.74 / (blank) / MD-2260.R00.00
$ws.Cells.Item(97, 1).Value = "Waiver"
$ws.Cells.Item(97, 2).Value = "MD"
$ws.Cells.Item(97, 3).Value = "1915(c)"
$ws.Cells.Item(97, 4).Value = "Amendment"
$ws.Cells.Item(97, 5).Value = "MD-2260.R00.74"
$ws.Cells.Item(97, 7).Value = "MD-2260.R00.00"
